$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 0.1473886666666666
$ws.Range("H2").Value = 0.4421659999999999
$ws.Range("O2").Value = 0.04602678343450817
$ws.Range("P2").Value = 0.06748706880158217
$ws.Range("Q2").Value = 0.004184511635333332
$ws.Range("R2").Value = 0.03766060471799999
$ws.Range("S2").Value = 0.04602678343450817
$ws.Range("T2").Value = 0.06748706880158217

# Row 3 updates
$ws.Range("G3").Value = 0.1473886666666666
$ws.Range("H3").Value = 0.4421659999999999
$ws.Range("M3").Value = 0.5884455
$ws.Range("N3").Value = 1.176891
$ws.Range("O3").Value = 0.9539732165654917
$ws.Range("P3").Value = 0.9325129311984178
$ws.Range("Q3").Value = 0.08673019765099997
$ws.Range("R3").Value = 0.5203811859059999
$ws.Range("S3").Value = 0.9539732165654917
$ws.Range("T3").Value = 0.9325129311984178
